# Update the mobile sample-metadata sheet with the latest Cocci batch
# (Result ID A1525801-A1525812 / lot 1405010 / run 20220913-Cocci-110317),
# replacing the previous batch's values in columns A, B, E, J, M, N, R, S, T, X
# for data rows 2-13. Columns that are unchanged (C, D, F, G, H, I, K, L, O, P,
# Q, U, V, W, Y) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$batchLot          = "'1405010"   # leading apostrophe: keep as text, not a number
$batchRun          = "20220913-Cocci-110317"
$batchAssay        = "A02804"
$batchComplex      = "TestComplexSite_20220913"
$batchFarm         = "TestFarm1_20220913"
$batchResultDate   = "'09/13/2022"  # leading apostrophe: keep as text, not a date
$batchResultTime   = "6:03 AM"
$batchCartridgeId  = "CartridgeCocci0317"
$batchCollectDate  = "'10/25/2022"  # leading apostrophe: keep as text, not a date

$resultIds = @(
    "A1525801",
    "A1525802",
    "A1525803",
    "A1525804",
    "A1525805",
    "A1525806",
    "A1525807",
    "A1525808",
    "A1525809",
    "A1525810",
    "A1525811",
    "A1525812"
)

for ($i = 0; $i -lt $resultIds.Length; $i++) {
    $row = 2 + $i

    $ws.Cells.Item($row, 1).Value  = $resultIds[$i]    # A - Result ID
    $ws.Cells.Item($row, 2).Value  = $batchLot          # B - Collection Site ID
    $ws.Cells.Item($row, 5).Value  = $batchRun          # E - Lab Sample ID
    $ws.Cells.Item($row, 10).Value = $batchAssay        # J - Flock ID
    $ws.Cells.Item($row, 13).Value = $batchComplex      # M - Complex
    $ws.Cells.Item($row, 14).Value = $batchFarm         # N - Farm
    $ws.Cells.Item($row, 18).Value = $batchResultDate   # R - Result Date
    $ws.Cells.Item($row, 19).Value = $batchResultTime   # S - Result Time
    $ws.Cells.Item($row, 20).Value = $batchCartridgeId  # T - Cartridge ID
    $ws.Cells.Item($row, 24).Value = $batchCollectDate  # X - Collection Date
}
